$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.205810785293579
$ws.Range("B1").Value = 2.714761018753052
$ws.Range("C1").Value = 2.005010366439819
$ws.Range("D1").Value = 1.866878986358643
$ws.Range("E1").Value = 1.724422693252563
